# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Force text storage so numeric-looking strings (e.g. "252.04", "4.50")
    # keep their exact formatting instead of being coerced to a Number,
    # then drop the temporary number-format override so the cell style
    # stays identical to its original (unstyled) state.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextCell $ws.Range('D2') '37.161.91'
$ws.Range('E2').Value = '  +0.15%  '
Set-TextCell $ws.Range('D3') '2.071.10'
$ws.Range('E3').Value = '  -1.03%  '
$ws.Range('E4').Value = '  +0.08%  '
Set-TextCell $ws.Range('D5') '252.04'
$ws.Range('E5').Value = '  +0.83%  '
Set-TextCell $ws.Range('D6') '0.675'
$ws.Range('E6').Value = '  +3.55%  '
Set-TextCell $ws.Range('D7') '61.24'
$ws.Range('E7').Value = '  +20.15%  '
$ws.Range('E8').Value = '  +0.07%  '
Set-TextCell $ws.Range('D9') '61.86'
$ws.Range('E9').Value = '  +2.42%  '
$ws.Range('E10').Value = '  +4.05%  '
Set-TextCell $ws.Range('D11') '0.0809'
$ws.Range('E11').Value = '  +9.48%  '
$ws.Range('E12').Value = '  +2.54%  '
Set-TextCell $ws.Range('D13') '15.94'
$ws.Range('E13').Value = '  +4.21%  '
Set-TextCell $ws.Range('D14') '2.376.64'
$ws.Range('E14').Value = '  +0.05%  '
Set-TextCell $ws.Range('D15') '0.819'
$ws.Range('E15').Value = '  -1.08%  '
Set-TextCell $ws.Range('D16') '5.46'
$ws.Range('E16').Value = '  +7.76%  '
Set-TextCell $ws.Range('D17') '2.078.32'
$ws.Range('E17').Value = '  -0.64%  '
Set-TextCell $ws.Range('D18') '37.130.00'
$ws.Range('E18').Value = '  +0.50%  '
Set-TextCell $ws.Range('D19') '74.68'
$ws.Range('E19').Value = '  +3.66%  '
Set-TextCell $ws.Range('D20') '15.36'
$ws.Range('E20').Value = '  +15.48%  '
Set-TextCell $ws.Range('D21') '0.0₃0925'
$ws.Range('E21').Value = '  +12.63%  '
Set-TextCell $ws.Range('D22') '5.50'
$ws.Range('E22').Value = '  +5.60%  '
Set-TextCell $ws.Range('D23') '239.72'
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('E25').Value = '  -0.09%  '
Set-TextCell $ws.Range('D26') '171.16'
$ws.Range('E26').Value = '  +1.06%  '
Set-TextCell $ws.Range('D27') '2.19'
$ws.Range('E27').Value = '  +9.01%  '
Set-TextCell $ws.Range('D28') '9.26'
$ws.Range('E28').Value = '  +1.25%  '
Set-TextCell $ws.Range('D29') '20.33'
$ws.Range('E29').Value = '  -1.79%  '
$ws.Range('E30').Value = '  +3.24%  '
Set-TextCell $ws.Range('D31') '4.75'
$ws.Range('E31').Value = '  +5.70%  '
$ws.Range('E32').Value = '  +3.44%  '
$ws.Range('E33').Value = '  +5.05%  '
Set-TextCell $ws.Range('D34') '4.40'
$ws.Range('E34').Value = '  +7.62%  '
Set-TextCell $ws.Range('D35') '0.0897'
$ws.Range('E35').Value = '  -1.47%  '
$ws.Range('E36').Value = '  +0.09%  '
Set-TextCell $ws.Range('D37') '2.29'
$ws.Range('E37').Value = '  -0.75%  '
$ws.Range('E38').Value = '  -2.59%  '
$ws.Range('E39').Value = '  +23.65%  '
$ws.Range('E40').Value = '  +2.39%  '
Set-TextCell $ws.Range('D41') '18.13'
$ws.Range('E41').Value = '  +1.85%  '
$ws.Range('E42').Value = '  +0.58%  '
$ws.Range('B43').Value = 'FTXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell $ws.Range('D43') '4.50'
$ws.Range('E43').Value = '  +28.49%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell $ws.Range('D44') '1.16'
$ws.Range('E44').Value = '  +0.42%  '
Set-TextCell $ws.Range('D45') '98.43'
$ws.Range('E45').Value = '  +0.20%  '
$ws.Range('E46').Value = '  +1.35%  '
Set-TextCell $ws.Range('D47') '4.67'
$ws.Range('E47').Value = '  +15.93%  '
$ws.Range('E48').Value = '  +10.02%  '
Set-TextCell $ws.Range('D49') '1.306.44'
$ws.Range('E49').Value = '  -0.45%  '
$ws.Range('E51').Value = '  +0.38%  '
